$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell updates from the cryptos list refresh.
# Cells whose new value would be auto-parsed as a pure number by Excel
# are forced to text (matching the original inline-string/text storage)
# by setting NumberFormat to "@" (Text) before assignment.

$ws.Range("D2").Value2 = '26.089.71'
$ws.Range("E2").Value2 = '  -0.66%  '
$ws.Range("D3").Value2 = '1.660.75'
$ws.Range("E3").Value2 = '  -1.40%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value2 = '1.003'
$ws.Range("E4").Value2 = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = '207.59'
$ws.Range("E5").Value2 = '  -2.42%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = '0.5172'
$ws.Range("E6").Value2 = '  -2.99%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = '1.003'
$ws.Range("E7").Value2 = '  -0.10%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = '0.2583'
$ws.Range("E8").Value2 = '  -4.20%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = '0.06299'
$ws.Range("E9").Value2 = '  -0.56%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = '20.98'
$ws.Range("E10").Value2 = '  -2.20%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = '0.07539'
$ws.Range("E11").Value2 = '  -0.06%  '
$ws.Range("D12").Value2 = '1.661.64'
$ws.Range("E12").Value2 = '  -1.46%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = '4.405'
$ws.Range("E13").Value2 = '  -2.26%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = '0.5369'
$ws.Range("E14").Value2 = '  -5.90%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = '66.09'
$ws.Range("E15").Value2 = '  -1.06%  '
$ws.Range("D16").Value2 = '0.0₅7926'
$ws.Range("E16").Value2 = '  -3.37%  '
$ws.Range("D17").Value2 = '26.120.41'
$ws.Range("E17").Value2 = '  -0.63%  '
$ws.Range("E18").Value2 = '  -0.09%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = '4.689'
$ws.Range("E19").Value2 = '  -3.88%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = '187.89'
$ws.Range("E20").Value2 = '  -1.15%  '
$ws.Range("E21").Value2 = '  -4.24%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = '6.176'
$ws.Range("E22").Value2 = '  -1.13%  '
$ws.Range("E23").Value2 = '  +0.07%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = '148.50'
$ws.Range("E24").Value2 = '  +0.34%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = '0.1214'
$ws.Range("E25").Value2 = '  -4.42%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = '7.384'
$ws.Range("E26").Value2 = '  -3.56%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = '15.58'
$ws.Range("E27").Value2 = '  -2.73%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = '1.382'
$ws.Range("E28").Value2 = '  +3.00%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = '0.06146'
$ws.Range("E29").Value2 = '  -5.11%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = '1.260'
$ws.Range("E30").Value2 = '  -2.20%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = '3.465'
$ws.Range("E31").Value2 = '  -2.33%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = '3.392'
$ws.Range("E32").Value2 = '  -3.05%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = '1.625'
$ws.Range("E33").Value2 = '  -2.33%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = '0.9863'
$ws.Range("E34").Value2 = '  -2.98%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = '2.388'
$ws.Range("E35").Value2 = '  -1.03%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = '2.750'
$ws.Range("E36").Value2 = '  +0.89%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = '0.5861'
$ws.Range("E37").Value2 = '  -4.58%  '
$ws.Range("D38").Value2 = '1.102.73'
$ws.Range("E38").Value2 = '  -0.50%  '
$ws.Range("E39").Value2 = '  -2.16%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = '5.976'
$ws.Range("E40").Value2 = '  -3.44%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = '0.8459'
$ws.Range("E41").Value2 = '  -2.70%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = '1.003'
$ws.Range("E42").Value2 = '  -0.52%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = '99.85'
$ws.Range("E43").Value2 = '  -0.45%  '
$ws.Range("D44").Value2 = '1.813.74'
$ws.Range("E44").Value2 = '  -1.04%  '
$ws.Range("D45").Value2 = '0.0₈106'
$ws.Range("E45").Value2 = '  -2.28%  '
$ws.Range("B46").Value2 = 'Frax'
$ws.Range("C46").Value2 = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = '1.002'
$ws.Range("E46").Value2 = '  -0.62%  '
$ws.Range("B47").Value2 = 'Aave'
$ws.Range("C47").Value2 = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = '54.88'
$ws.Range("E47").Value2 = '  -3.94%  '
$ws.Range("B48").Value2 = 'EnergySwap'
$ws.Range("C48").Value2 = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = '8.012'
$ws.Range("E48").Value2 = '  -0.55%  '
$ws.Range("B49").Value2 = 'Cronos'
$ws.Range("C49").Value2 = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = '0.05245'
$ws.Range("E49").Value2 = '  -0.46%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = '0.4244'
$ws.Range("E50").Value2 = '  -0.73%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = '5.858'
$ws.Range("E51").Value2 = '  -2.20%  '
